# Weekly fruit/vegetable update: insert a new week of data (rows 366-367)
# and push the previously-existing rows 366-386 down to 368-388.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 366-367; everything below shifts down by 2
# (old row 366 -> 368, old row 386 -> 388).
$ws.Range("A366:A367").EntireRow.Insert()

# New row 366: Primera
$ws.Cells.Item(366, 1).Value = 11
$ws.Cells.Item(366, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(366, 3).Value = 'Bíobío'
$ws.Cells.Item(366, 4).Value = 45008
$ws.Cells.Item(366, 5).Value = 8
$ws.Cells.Item(366, 6).Value = 100112009
$ws.Cells.Item(366, 7).Value = 'Acelga'
$ws.Cells.Item(366, 8).Value = 'Sin especificar'
$ws.Cells.Item(366, 9).Value = 'Primera'
$ws.Cells.Item(366, 10).Value = 200
$ws.Cells.Item(366, 11).Value = 700
$ws.Cells.Item(366, 12).Value = 800
$ws.Cells.Item(366, 13).Value = 750
$ws.Cells.Item(366, 14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(366, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(366, 16).Value = 750
$ws.Cells.Item(366, 17).Value = 1
$ws.Cells.Item(366, 18).Value = 'Hortaliza'

# New row 367: Segunda
$ws.Cells.Item(367, 1).Value = 11
$ws.Cells.Item(367, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(367, 3).Value = 'Bíobío'
$ws.Cells.Item(367, 4).Value = 45008
$ws.Cells.Item(367, 5).Value = 8
$ws.Cells.Item(367, 6).Value = 100112009
$ws.Cells.Item(367, 7).Value = 'Acelga'
$ws.Cells.Item(367, 8).Value = 'Sin especificar'
$ws.Cells.Item(367, 9).Value = 'Segunda'
$ws.Cells.Item(367, 10).Value = 100
$ws.Cells.Item(367, 11).Value = 600
$ws.Cells.Item(367, 12).Value = 600
$ws.Cells.Item(367, 13).Value = 600
$ws.Cells.Item(367, 14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(367, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(367, 16).Value = 600
$ws.Cells.Item(367, 17).Value = 1
$ws.Cells.Item(367, 18).Value = 'Hortaliza'
